$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 3 blank rows before the old blank-separator row (row 8),
#    pushing the blank separator row (old 8) down to row 11 and the
#    header row (old 9) down to row 12 -- preserving their styles.
# ------------------------------------------------------------------
$ws.Rows("8:10").Insert()

# ------------------------------------------------------------------
# 2. Update the "Howto" steps in column D (rows 3-10).
#    NOTE: cells are written in the same order as the original author
#    entered them, so new shared-string entries land at the same
#    indices as the canonical file (D3, D5, D4, D6, D7, D10, D8, D9).
# ------------------------------------------------------------------
$ws.Range("D3").Value = "1. Put a copy of this file in a folder."
$ws.Range("D5").Value = "3. Assign increasing AnonIDs to your Cases, for example by typing in the one on the first row (eg ""MYPROJ-001"") and then dragging down to number the rest sequentially."
$ws.Range("D4").Value = "2. List Cases to anonymize below. Use the same .zip filenames or folder names as you get when exporting cases from the PACS, preferably in random order."
$ws.Range("D6").Value = "4. Export cases and put in same folder. If using folders, these should have BLOCK_STAIN subfolders cotaining one .svs or .ndpi file each. Don't fill up more than half your available storage."
$ws.Range("D7").Value = "5. Run aida-pat-anonexcel.py on this file to check for mistakes, anonymize slides, and update this sheet to match."

$ws.Range("D10").Value = "8. Export to research system. Delete exported cases and anonymized images. Repeat from 4 until all cases have been processed."

# Row 8 (new): rich text - "6. " + bold "Your output data is now Pseudonymous" + rest
$cell8 = $ws.Range("D8")
$cell8.Value = "6. Your output data is now Pseudonymous because keys still exist that connect AnonIDs to persons. Take this moment to verify that everything went as expected."
$cell8.Characters(1, 3).Font.Size = 11
$cell8.Characters(4, 36).Font.Bold = $true
$cell8.Characters(40, 119).Font.Size = 11

# Row 9 (new): rich text with several bold runs
$cell9 = $ws.Range("D9")
$cell9.Value = "7. To make your data Anonymous: Delete all keys associating AnonIDs to persons, including the Case and OrigFile cells below and any other identifiers that may exist. Obviously, Study parameters may not contain identifiers."
$cell9.Characters(1, 3).Font.Size = 11
$cell9.Characters(4, 27).Font.Bold = $true
$cell9.Characters(31, 64).Font.Size = 11
$cell9.Characters(95, 4).Font.Bold = $true
$cell9.Characters(99, 5).Font.Size = 11
$cell9.Characters(104, 8).Font.Bold = $true
$cell9.Characters(112, 66).Font.Size = 11
$cell9.Characters(178, 16).Font.Bold = $true
$cell9.Characters(194, 29).Font.Size = 11

# ------------------------------------------------------------------
# 3. Rename the "Person" header (column B, now row 12) to "Case".
# ------------------------------------------------------------------
$ws.Range("B12").Value = "Case"

# ------------------------------------------------------------------
# 4. Update the window view position + active selection.
# ------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 4845
$win.Top = 4335

$ws.Range("B12").Select()
